# Append a new Q&A row (row 12) to the query_responses sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "What is 1 + 1?"
# Prefix with an apostrophe so the numeric-looking answer "2." is stored
# as literal text rather than being coerced to the number 2.
$ws.Range("B12").Value = "'2."
